# Add 25 new rows (721-745) to Sheet1, mirroring the pattern of the
# existing "remote api config" rows already present at 717-720.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$imageFileName = "0a5dbcdc-2e44-4579-a576-c93d5ee55485.png"

# Rows 721 through 744: A = 0.0, B = image file name (shared string index 1)
for ($r = 721; $r -le 744; $r++) {
    $ws.Cells.Item($r, 1).Value = 0.0
    $ws.Cells.Item($r, 2).Value = $imageFileName
}

# Row 745: A = 126.0, B = image file name (shared string index 1)
$ws.Cells.Item(745, 1).Value = 126.0
$ws.Cells.Item(745, 2).Value = $imageFileName
